{"js": "// Office.js (Word JavaScript API) script.\n// Applies the resume \"job description enhancement\" edit described in the\n// commit message / XML diff: a series of in-place bullet/heading text\n// replacements plus a handful of newly-added bullet paragraphs.\n//\n// The document has one run/one <w:t> per paragraph, so we operate purely\n// in terms of `context.document.body.paragraphs` items, addressed by their\n// (stable, 0-based) index in document order. Using indices avoids any\n// ambiguity from the \"Political Research and Data Analysis\" heading text\n// that appears twice in the original document but must become two\n// different things.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Sanity-check a handful of paragraphs' text before mutating, so that if\n// the document shape ever changes this script fails loudly instead of\n// silently corrupting unrelated content.\nconst checkIdx = [9, 10, 11, 12, 13, 14, 17, 18, 19, 20, 21, 33, 36, 37, 38, 39, 53, 54, 55, 56, 57, 59, 60, 61, 62, 63];\ncheckIdx.forEach(i => paragraphs.items[i].load(\"text\"));\nawait context.sync();\n\nconst expected = {\n  9:  \"\\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations\",\n  10: \"\\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics\",\n  11: \"\\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\",\n  12: \"\\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\",\n  13: \"\\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\",\n  14: \"\\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\",\n  17: \"\\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\",\n  18: \"\\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\",\n  19: \"\\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\",\n  20: \"\\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products\",\n  21: \"\\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\",\n  33: \"\\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research\",\n  36: \"\\u2022 Managed critical research operations for political campaigns\",\n  37: \"\\u2022 Conducted comprehensive polling and demographic analysis\",\n  38: \"\\u2022 Developed strategic recommendations based on data analysis\",\n  39: \"\\u2022 Led research team in support of progressive political initiatives\",\n  53: \"Political Research and Data Analysis\",\n  54: \"\\u2022 Developed data analysis tools for political polling and research\",\n  55: \"\\u2022 Built statistical models for voter behavior analysis\",\n  56: \"\\u2022 Created data visualization tools for research presentations\",\n  57: \"\\u2022 Supported senior researchers with technical analysis and reporting\",\n  59: \"Political Field Operations and Data Management\",\n  60: \"\\u2022 Managed field operations for political campaigns and research projects\",\n  61: \"\\u2022 Developed data collection and management systems for field work\",\n  62: \"\\u2022 Trained field staff on data collection protocols and quality control\",\n  63: \"\\u2022 Analyzed field data to inform campaign strategy and research findings\",\n};\n\nfor (const i of checkIdx) {\n  const actual = paragraphs.items[i].text;\n  if (actual !== expected[i]) {\n    throw new Error(\n      \"Unexpected content at paragraph \" + i + \": \" + JSON.stringify(actual)\n    );\n  }\n}\n\n// Helper: replace the full text of a paragraph in place (keeps the\n// paragraph's own formatting / style, matches the diff's \"modified line\"\n// semantics).\nfunction setText(idx, newText) {\n  paragraphs.items[idx].insertText(newText, Word.InsertLocation.replace);\n}\n\n// --- PARTNER - Siege Analytics bullets ---\nsetText(9,  \"\\u2022 Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\");\nsetText(10, \"\\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis\");\nsetText(11, \"\\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\");\nsetText(12, \"\\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\");\nsetText(13, \"\\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\");\nsetText(14, \"\\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\");\n\n// --- DATA PRODUCTS MANAGER - Helm/Murmuration bullets ---\nsetText(17, \"\\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\");\nsetText(18, \"\\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\");\nsetText(19, \"\\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\");\nsetText(20, \"\\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products\");\nsetText(21, \"\\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\");\n\n// --- SENIOR ANALYST - Myers Research bullet ---\nsetText(33, \"\\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions\");\n\n// --- RESEARCH DIRECTOR - Progressive Change Campaign Committee bullets ---\nsetText(36, \"\\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\");\nsetText(37, \"\\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\");\nsetText(38, \"\\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\");\nsetText(39, \"\\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs\");\n// New bullet added after paragraph 39.\nparagraphs.items[39].insertParagraph(\n  \"\\u2022 Managed comprehensive research operations for progressive political initiatives and candidates\",\n  Word.InsertLocation.after\n);\n\n// --- PROGRAMMER - Lake Research Partners heading + bullets ---\nsetText(53, \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\");\nsetText(54, \"\\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\");\nsetText(55, \"\\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute\");\nsetText(56, \"\\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions\");\n// Two new bullets added after paragraph 56 (insert in reverse order so the\n// final document order is \"Conducted statistical...\" then \"Pioneered...\").\nparagraphs.items[56].insertParagraph(\n  \"\\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps\",\n  Word.InsertLocation.after\n);\nparagraphs.items[56].insertParagraph(\n  \"\\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\",\n  Word.InsertLocation.after\n);\nsetText(57, \"\\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\");\n\n// --- FIELD DIRECTOR - The Feldman Group heading + bullets ---\nsetText(59, \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\");\nsetText(60, \"\\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions\");\nsetText(61, \"\\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm\");\nsetText(62, \"\\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings\");\nsetText(63, \"\\u2022 Created custom reports and data visualizations based on specific client requirements\");\n// Two new bullets added after paragraph 63 (insert in reverse order, as above).\nparagraphs.items[63].insertParagraph(\n  \"\\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\",\n  Word.InsertLocation.after\n);\nparagraphs.items[63].insertParagraph(\n  \"\\u2022 Introduced mapping and geospatial analysis into standard reporting procedures\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the resume \"job description enhancement\" edit described in the\n# commit message / XML diff: a series of in-place bullet/heading text\n# replacements plus a handful of newly-added bullet paragraphs.\n#\n# The document has one run/one <w:t> per paragraph, so we operate purely\n# via $d.Paragraphs.Item(N) (1-based indices, matching Word's own COM\n# numbering). We address paragraphs by their stable position in document\n# order rather than by searching for text, because the heading \"Political\n# Research and Data Analysis\" appears twice in the original document and\n# must turn into two different things depending on which job it belongs to.\n#\n# NOTE: $d.Paragraphs.Item(N) is 1-based, i.e. it equals the 0-based\n# paragraph index (as you'd get from Office.js) plus 1. All string\n# literals below use single quotes (with '' for an embedded apostrophe)\n# so the literal $ in one of the new bullets doesn't need escaping.\n\n$d = $word.ActiveDocument\n\n# --- Sanity-check a handful of paragraphs before mutating -------------\n$expected = @{\n  10 = '\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations'\n  11 = '\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics'\n  12 = '\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets'\n  13 = '\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering'\n  14 = '\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications'\n  15 = '\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices'\n  18 = '\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES'\n  19 = '\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions'\n  20 = '\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI'\n  21 = '\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company''s distinguishing products'\n  22 = '\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices'\n  34 = '\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research'\n  37 = '\u2022 Managed critical research operations for political campaigns'\n  38 = '\u2022 Conducted comprehensive polling and demographic analysis'\n  39 = '\u2022 Developed strategic recommendations based on data analysis'\n  40 = '\u2022 Led research team in support of progressive political initiatives'\n  54 = 'Political Research and Data Analysis'\n  55 = '\u2022 Developed data analysis tools for political polling and research'\n  56 = '\u2022 Built statistical models for voter behavior analysis'\n  57 = '\u2022 Created data visualization tools for research presentations'\n  58 = '\u2022 Supported senior researchers with technical analysis and reporting'\n  60 = 'Political Field Operations and Data Management'\n  61 = '\u2022 Managed field operations for political campaigns and research projects'\n  62 = '\u2022 Developed data collection and management systems for field work'\n  63 = '\u2022 Trained field staff on data collection protocols and quality control'\n  64 = '\u2022 Analyzed field data to inform campaign strategy and research findings'\n}\n\nforeach ($idx in $expected.Keys) {\n  $actual = $d.Paragraphs.Item($idx).Range.Text.TrimEnd([char]13)\n  if ($actual -ne $expected[$idx]) {\n    throw \"Unexpected content at paragraph $idx : $actual\"\n  }\n}\n\n# --- PARTNER - Siege Analytics bullets ---\n$d.Paragraphs.Item(10).Range.Text = '\u2022 Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions'\n$d.Paragraphs.Item(11).Range.Text = '\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis'\n$d.Paragraphs.Item(12).Range.Text = '\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets'\n$d.Paragraphs.Item(13).Range.Text = '\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering'\n$d.Paragraphs.Item(14).Range.Text = '\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications'\n$d.Paragraphs.Item(15).Range.Text = '\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices'\n\n# --- DATA PRODUCTS MANAGER - Helm/Murmuration bullets ---\n$d.Paragraphs.Item(18).Range.Text = '\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES'\n$d.Paragraphs.Item(19).Range.Text = '\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions'\n$d.Paragraphs.Item(20).Range.Text = '\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI'\n$d.Paragraphs.Item(21).Range.Text = '\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company''s distinguishing products'\n$d.Paragraphs.Item(22).Range.Text = '\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices'\n\n# --- SENIOR ANALYST - Myers Research bullet ---\n$d.Paragraphs.Item(34).Range.Text = '\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions'\n\n# --- RESEARCH DIRECTOR - Progressive Change Campaign Committee bullets ---\n$d.Paragraphs.Item(37).Range.Text = '\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls'\n$d.Paragraphs.Item(38).Range.Text = '\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren'\n$d.Paragraphs.Item(39).Range.Text = '\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver'\n$d.Paragraphs.Item(40).Range.Text = '\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs'\n# New bullet added after paragraph 40 (the \"Designed survey deployment...\" one).\n$d.Paragraphs.Item(40).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(41).Range.Text = '\u2022 Managed comprehensive research operations for progressive political initiatives and candidates'\n\n# --- PROGRAMMER - Lake Research Partners heading + bullets ---\n# (shifted by +1 from the original numbering because of the insertion above)\n$d.Paragraphs.Item(55).Range.Text = 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns'\n$d.Paragraphs.Item(56).Range.Text = '\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party'\n$d.Paragraphs.Item(57).Range.Text = '\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute'\n$d.Paragraphs.Item(58).Range.Text = '\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions'\n# Two new bullets added after paragraph 58 (inserted forward in order).\n$d.Paragraphs.Item(58).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(59).Range.Text = '\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle'\n$d.Paragraphs.Item(59).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(60).Range.Text = '\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps'\n$d.Paragraphs.Item(61).Range.Text = '\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding'\n\n# --- FIELD DIRECTOR - The Feldman Group heading + bullets ---\n# (shifted by +3 total now from the two insertions above)\n$d.Paragraphs.Item(63).Range.Text = 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns'\n$d.Paragraphs.Item(64).Range.Text = '\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions'\n$d.Paragraphs.Item(65).Range.Text = '\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm'\n$d.Paragraphs.Item(66).Range.Text = '\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings'\n$d.Paragraphs.Item(67).Range.Text = '\u2022 Created custom reports and data visualizations based on specific client requirements'\n# Two new bullets added after paragraph 67 (inserted forward in order).\n$d.Paragraphs.Item(67).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(68).Range.Text = '\u2022 Introduced mapping and geospatial analysis into standard reporting procedures'\n$d.Paragraphs.Item(68).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(69).Range.Text = '\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL'\n"}
